# Scheduled market-price refresh for Pandaemonium Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per
# crafting-job leve row, sourced from the latest market-board snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 52: Your Courtesy Wake-up Call (Smelling Salts)
$ws.Range("H52").Value = 91667.27
$ws.Range("J52").Value = 100784
$ws.Range("L52").Value = 302352
$ws.Range("N52").Value = -302672

# Row 58: A Matter of Vital Importance (Mega-Potion of Vitality)
$ws.Range("H58").Value = 37.5
$ws.Range("I58").Value = 37.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 112.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 37.5
$ws.Range("N58").ClearContents()

# Row 112: Making Ends Meet (Superior Spiritbond Potion)
$ws.Range("H112").Value = 1551.0358
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 1593.6666
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 4780.9998
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -6996.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 88: The Mast Chance (Adamantite Rivets)
$ws.Range("H88").Value = 4931
$ws.Range("I88").Value = 18370.666
$ws.Range("J88").Value = 1829.5385
$ws.Range("K88").Value = 18370.666
$ws.Range("L88").Value = 1829.5385
$ws.Range("M88").Value = -17964.666
$ws.Range("N88").Value = -2641.5385

# Row 91: The Rose and the Riveter (L) (Adamantite Rivets)
$ws.Range("H91").Value = 4931
$ws.Range("I91").Value = 18370.666
$ws.Range("J91").Value = 1829.5385
$ws.Range("K91").Value = 18370.666
$ws.Range("L91").Value = 1829.5385
$ws.Range("M91").Value = -16966.666
$ws.Range("N91").Value = -4637.538500000001

# Row 122: Haste for High Durium (High Durium Nugget)
$ws.Range("H122").Value = 7354906.5
$ws.Range("I122").Value = 1958
$ws.Range("J122").Value = 41668668
$ws.Range("K122").Value = 5874
$ws.Range("L122").Value = 125006004
$ws.Range("M122").Value = -3424
$ws.Range("N122").Value = -125010904

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin (Adamantite Nugget)
$ws.Range("H86").Value = 2200.2917
$ws.Range("I86").Value = 2300
$ws.Range("J86").Value = 1502.3334
$ws.Range("K86").Value = 2300
$ws.Range("L86").Value = 1502.3334
$ws.Range("M86").Value = -1177
$ws.Range("N86").Value = -3748.3334

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) (Adamantite Nugget)
$ws.Range("H89").Value = 2200.2917
$ws.Range("I89").Value = 2300
$ws.Range("J89").Value = 1502.3334
$ws.Range("K89").Value = 11500
$ws.Range("L89").Value = 7511.666999999999
$ws.Range("M89").Value = -5884
$ws.Range("N89").Value = -18743.667

$ws = $wb.Worksheets.Item("CRP")
# Row 51: Greenstone for Greenhorns (Jade Crook)
$ws.Range("H51").Value = 23471
$ws.Range("J51").Value = 23471
$ws.Range("L51").Value = 23471
$ws.Range("N51").Value = -24943

# Row 61: Incant Now, Think Later (Jade Crook)
$ws.Range("H61").Value = 23471
$ws.Range("J61").Value = 23471
$ws.Range("L61").Value = 23471
$ws.Range("N61").Value = -24167

# Row 99: O Pine (Pine Lumber)
$ws.Range("H99").Value = 1800
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# Row 107: Built to Last (White Oak Lumber)
$ws.Range("H107").Value = 1536.7273
$ws.Range("I107").Value = 1511.5555
$ws.Range("J107").Value = 1650
$ws.Range("K107").Value = 1511.5555
$ws.Range("L107").Value = 1650
$ws.Range("M107").Value = 408.4445000000001
$ws.Range("N107").Value = -5490

# Row 122: Timber of Tenkonto (Horse Chestnut Lumber)
$ws.Range("H122").Value = 11124.381
$ws.Range("I122").Value = 4833.25
$ws.Range("J122").Value = 19512.555
$ws.Range("K122").Value = 14499.75
$ws.Range("L122").Value = 58537.665
$ws.Range("M122").Value = -12049.75
$ws.Range("N122").Value = -63437.665

# Row 126: A Better Conductor (Red Pine Lumber)
$ws.Range("H126").Value = 1800
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 2005.1316
$ws.Range("I132").Value = 1814.8438
$ws.Range("K132").Value = 5444.5314
$ws.Range("M132").Value = -2914.5314

$ws = $wb.Worksheets.Item("CUL")
# Row 25: Flakes for Friends (Apple Tart)
$ws.Range("H25").Value = 4466.067
$ws.Range("I25").Value = 995.5
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 2986.5
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = -2817.5
$ws.Range("N25").Value = -15338

# Row 30: Picnic Panic (Apple Tart)
$ws.Range("H30").Value = 4466.067
$ws.Range("I30").Value = 995.5
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 2986.5
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = -2884.5
$ws.Range("N30").Value = -15204

# Row 98: Sweet Kiss of Death (Rice Vinegar)
$ws.Range("H98").Value = 423.21622
$ws.Range("I98").Value = 330.92593
$ws.Range("J98").Value = 672.4
$ws.Range("K98").Value = 992.77779
$ws.Range("L98").Value = 2017.2
$ws.Range("M98").Value = 505.22221
$ws.Range("N98").Value = -5013.2

# Row 105: Fish Box (Chirashi-zushi)
$ws.Range("H105").Value = 6800
$ws.Range("J105").Value = 6800
$ws.Range("L105").Value = 20400
$ws.Range("N105").Value = -25642

# Row 121: A Cookie for Your Troubles (Coffee Biscuit)
$ws.Range("H121").Value = 645
$ws.Range("I121").Value = 645
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 1935
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -625
$ws.Range("N121").ClearContents()

# Row 122: Salt of the North (Northern Sea Salt)
$ws.Range("H122").Value = 963.94116
$ws.Range("I122").Value = 618.9
$ws.Range("J122").Value = 1456.8572
$ws.Range("K122").Value = 5570.099999999999
$ws.Range("L122").Value = 13111.7148
$ws.Range("M122").Value = -3120.099999999999
$ws.Range("N122").Value = -18011.7148

# Row 126: Imperial Palate (Glory Be Soup)
$ws.Range("H126").Value = 1562.3334
$ws.Range("I126").Value = 1050.6428
$ws.Range("J126").Value = 2585.7144
$ws.Range("K126").Value = 3151.9284
$ws.Range("L126").Value = 7757.1432
$ws.Range("M126").Value = 1788.0716
$ws.Range("N126").Value = -17637.1432

# Row 129: Comfort Food (Yakow Moussaka)
$ws.Range("H129").Value = 2449.8235
$ws.Range("I129").Value = 3268.625
$ws.Range("J129").Value = 1722
$ws.Range("K129").Value = 9805.875
$ws.Range("L129").Value = 5166
$ws.Range("M129").Value = -4805.875
$ws.Range("N129").Value = -15166

$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me (Bone Hora)
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 102: Put the Metal to the Peddle (Durium Ingot)
$ws.Range("H102").Value = 3887.476
$ws.Range("I102").Value = 3445.7812
$ws.Range("K102").Value = 3445.7812
$ws.Range("M102").Value = -1823.7812

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs (Aldgoat Leather)
$ws.Range("H22").Value = 672.63635
$ws.Range("I22").Value = 866
$ws.Range("J22").Value = 600.125
$ws.Range("K22").Value = 866
$ws.Range("L22").Value = 600.125
$ws.Range("M22").Value = -571
$ws.Range("N22").Value = -1190.125

# Row 27: Fire and Hide (Aldgoat Leather)
$ws.Range("H27").Value = 672.63635
$ws.Range("I27").Value = 866
$ws.Range("J27").Value = 600.125
$ws.Range("K27").Value = 866
$ws.Range("L27").Value = 600.125
$ws.Range("M27").Value = -759
$ws.Range("N27").Value = -814.125

$ws = $wb.Worksheets.Item("WVR")
# Row 112: Hair Do No Harm (Iridescent Hat of Healing)
$ws.Range("H112").Value = 79800
$ws.Range("J112").Value = 79800
$ws.Range("L112").Value = 79800
$ws.Range("N112").Value = -82754

# Row 115: Gloves Come in Handy (Pixie Cotton Sleeves of Crafting)
$ws.Range("H115").Value = 53950
$ws.Range("J115").Value = 53950
$ws.Range("L115").Value = 53950
$ws.Range("N115").Value = -57084

# Row 122: Heavy Armoire (Dark Hempen Cloth)
$ws.Range("H122").Value = 3504.325
$ws.Range("I122").Value = 2031.0646
$ws.Range("J122").Value = 8578.888999999999
$ws.Range("K122").Value = 6093.1938
$ws.Range("L122").Value = 25736.667
$ws.Range("M122").Value = -3643.1938
$ws.Range("N122").Value = -30636.667
